$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 304, shifting all subsequent rows (304-397) down
# to (305-398).
$ws.Rows.Item(304).Insert()

# Populate the newly inserted row 304 with the new data record.
$ws.Cells.Item(304, 1).Value = 5
$ws.Cells.Item(304, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(304, 3).Value = "Maule"
$ws.Cells.Item(304, 4).Value = 44588
$ws.Cells.Item(304, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(304, 5).Value = 7
$ws.Cells.Item(304, 6).Value = 100114001
$ws.Cells.Item(304, 7).Value = "Papa"
$ws.Cells.Item(304, 8).Value = "Patagonia"
$ws.Cells.Item(304, 9).Value = "1a nueva(o)"
$ws.Cells.Item(304, 10).Value = 1200
$ws.Cells.Item(304, 11).Value = 7000
$ws.Cells.Item(304, 12).Value = 7000
$ws.Cells.Item(304, 13).Value = 7000
$ws.Cells.Item(304, 14).Value = "`$/saco 25 kilos"
$ws.Cells.Item(304, 15).Value = "Región de Los Lagos"
$ws.Cells.Item(304, 16).Value = 280
$ws.Cells.Item(304, 17).Value = 25
$ws.Cells.Item(304, 18).Value = "Hortaliza"
